$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1) - reuse the same header style as the
# existing header cells (copy format from G1, then set the text).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H15 with 0 (numeric) to match the new "Save" column data.
$ws.Range("H2:H15").Value = 0
